$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 788 (shifts rows 788:829 down to 789:830)
$ws.Rows.Item(788).Insert()

# Populate the newly inserted row with the new data point.
# A788/B788 share the same date ("2026/02/08") and weekday ("日") as row
# 787, which already holds that text literally - copy it across instead of
# assigning a literal string so Excel doesn't reinterpret the slash-
# separated text as a date serial number.
$ws.Range("A787").Copy($ws.Range("A788"))
$ws.Range("B787").Copy($ws.Range("B788"))
$ws.Range("C788").Value = 8
$ws.Range("D788").Value = 201
